# Update the "Rules" workbook: change the greeting text in E8 from
# "Good Morning" to "GIT UPDATE", and leave the selection on that cell
# (mirrors what Excel records when a user edits a cell and the cursor
# stays there).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("E8").Value = "GIT UPDATE"
$ws.Range("E8").Select()
